# [Kadastro App] Yeni kayit eklendi: 3016
# Adds the new record row (row 75) to both the master "Kayitlar" sheet
# and the per-district "Erdemli" sheet, mirroring the author's edit.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "3016"
    B = "2025-09-11"
    C = "Erdemli"
    D = "1"
    E = "3B"
    F = "SERDAR ARSLAN (Tekniker), ÖZKAN AKBAŞ (Mühendis)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # All columns in this table are stored as plain text (even the
    # numeric-looking "Kayıt No" / "Parsel Sayısı" and the date string),
    # so use a leading apostrophe to force text entry, exactly like the
    # existing rows above it, instead of letting Excel auto-convert
    # "3016" -> number or "2025-09-11" -> date.
    $ws.Range("A75").Value = "'" + $newRow.A
    $ws.Range("B75").Value = "'" + $newRow.B
    $ws.Range("C75").Value = "'" + $newRow.C
    $ws.Range("D75").Value = "'" + $newRow.D
    $ws.Range("E75").Value = "'" + $newRow.E
    $ws.Range("F75").Value = "'" + $newRow.F
}
